# "generate trainer und segler" -- add two more tracked days (columns D and E)
# to the time-tracking sheet, extend the totals/hours formulas to cover the
# new columns (through F, which stays empty but keeps the formula pattern),
# and move the active selection to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date headers for the additional days, formatted like the existing C2 date.
$ws.Range("D2").Value = 42069
$ws.Range("D2").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("E2").Value = 42071
$ws.Range("E2").NumberFormat = $ws.Range("C2").NumberFormat

# New minutes-per-task entries for the added days.
$ws.Range("D3").Value = 15
$ws.Range("D6").Value = 225
$ws.Range("E6").Value = 60

# Extend the "total minutes" row (13) and "total hours" row (14) formulas
# across the new columns D, E and F (F being the next, still-empty column).
$ws.Range("D13").Formula = "=D3+D4+D5+D6"
$ws.Range("E13").Formula = "=E3+E4+E5+E6"
$ws.Range("F13").Formula = "=F3+F4+F5+F6"

$ws.Range("D14").Formula = "=D13/60"
$ws.Range("E14").Formula = "=E13/60"
$ws.Range("F14").Formula = "=F13/60"

# The author's cursor ends up on F2, the next empty date-header cell.
$ws.Range("F2").Select()
